# Update NATMI ligand-receptor pair metrics (Cp-Slc40a1) with refreshed TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 4.619095000000001
$ws.Range("H2").Value = 13.857285
$ws.Range("I2").Value = 0.09937071456472289
$ws.Range("J2").Value = 0.09937071456472289
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.270625
$ws.Range("N2").Value = 0.811875
$ws.Range("O2").Value = 0.00651929904432436
$ws.Range("P2").Value = 0.006519299044324359
$ws.Range("Q2").Value = 1.250042584375
$ws.Range("R2").Value = 11.250383259375
$ws.Range("S2").Value = 0.0006478274044956268
$ws.Range("T2").Value = 0.0006478274044956267

# Row 3
$ws.Range("G3").Value = 4.619095000000001
$ws.Range("H3").Value = 13.857285
$ws.Range("I3").Value = 0.09937071456472289
$ws.Range("J3").Value = 0.09937071456472289
$ws.Range("O3").Value = 0.2393219402230525
$ws.Range("P3").Value = 0.2393219402230525
$ws.Range("Q3").Value = 45.88877034479835
$ws.Range("R3").Value = 412.998933103185
$ws.Range("S3").Value = 0.02378159221098063
$ws.Range("T3").Value = 0.02378159221098063

# Row 4
$ws.Range("G4").Value = 4.619095000000001
$ws.Range("H4").Value = 13.857285
$ws.Range("I4").Value = 0.09937071456472289
$ws.Range("J4").Value = 0.09937071456472289
$ws.Range("M4").Value = 14.774121
$ws.Range("N4").Value = 44.322363
$ws.Range("O4").Value = 0.3559054518837227
$ws.Range("P4").Value = 0.3559054518837226
$ws.Range("Q4").Value = 68.24306844049501
$ws.Range("R4").Value = 614.1876159644549
$ws.Range("S4").Value = 0.03536657907116612
$ws.Range("T4").Value = 0.03536657907116612

# Row 5
$ws.Range("G5").Value = 4.619095000000001
$ws.Range("H5").Value = 13.857285
$ws.Range("I5").Value = 0.09937071456472289
$ws.Range("J5").Value = 0.09937071456472289
$ws.Range("M5").Value = 16.53203833333334
$ws.Range("N5").Value = 49.596115
$ws.Range("O5").Value = 0.3982533088489005
$ws.Range("P5").Value = 0.3982533088489005
$ws.Range("Q5").Value = 76.36305560530836
$ws.Range("R5").Value = 687.2675004477751
$ws.Range("S5").Value = 0.03957471587808052
$ws.Range("T5").Value = 0.03957471587808052

# Row 6
$ws.Range("I6").Value = 0.3438177451937012
$ws.Range("J6").Value = 0.3438177451937012
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.270625
$ws.Range("N6").Value = 0.811875
$ws.Range("O6").Value = 0.00651929904432436
$ws.Range("P6").Value = 0.006519299044324359
$ws.Range("Q6").Value = 4.325085359791667
$ws.Range("R6").Value = 38.925768238125
$ws.Range("S6").Value = 0.002241450697663053
$ws.Range("T6").Value = 0.002241450697663053

# Row 7
$ws.Range("I7").Value = 0.3438177451937012
$ws.Range("J7").Value = 0.3438177451937012
$ws.Range("O7").Value = 0.2393219402230525
$ws.Range("P7").Value = 0.2393219402230525
$ws.Range("S7").Value = 0.08228312986287167
$ws.Range("T7").Value = 0.08228312986287166

# Row 8
$ws.Range("I8").Value = 0.3438177451937012
$ws.Range("J8").Value = 0.3438177451937012
$ws.Range("M8").Value = 14.774121
$ws.Range("N8").Value = 44.322363
$ws.Range("O8").Value = 0.3559054518837227
$ws.Range("P8").Value = 0.3559054518837226
$ws.Range("Q8").Value = 236.117633037933
$ws.Range("R8").Value = 2125.058697341397
$ws.Range("S8").Value = 0.1223666099688069
$ws.Range("T8").Value = 0.1223666099688068

# Row 9
$ws.Range("I9").Value = 0.3438177451937012
$ws.Range("J9").Value = 0.3438177451937012
$ws.Range("M9").Value = 16.53203833333334
$ws.Range("N9").Value = 49.596115
$ws.Range("O9").Value = 0.3982533088489005
$ws.Range("P9").Value = 0.3982533088489005
$ws.Range("Q9").Value = 264.2123860065206
$ws.Range("R9").Value = 2377.911474058686
$ws.Range("S9").Value = 0.1369265546643597
$ws.Range("T9").Value = 0.1369265546643597

# Row 10
$ws.Range("G10").Value = 19.63122766666666
$ws.Range("H10").Value = 58.893683
$ws.Range("I10").Value = 0.4223271270713038
$ws.Range("J10").Value = 0.4223271270713038
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.270625
$ws.Range("N10").Value = 0.811875
$ws.Range("O10").Value = 0.00651929904432436
$ws.Range("P10").Value = 0.006519299044324359
$ws.Range("Q10").Value = 5.312700987291666
$ws.Range("R10").Value = 47.814308885625
$ws.Range("S10").Value = 0.002753276835908203
$ws.Range("T10").Value = 0.002753276835908203

# Row 11
$ws.Range("G11").Value = 19.63122766666666
$ws.Range("H11").Value = 58.893683
$ws.Range("I11").Value = 0.4223271270713038
$ws.Range("J11").Value = 0.4223271270713038
$ws.Range("O11").Value = 0.2393219402230525
$ws.Range("P11").Value = 0.2393219402230525
$ws.Range("Q11").Value = 195.0280082964559
$ws.Range("R11").Value = 1755.252074668103
$ws.Range("S11").Value = 0.1010721474595321
$ws.Range("T11").Value = 0.1010721474595321

# Row 12
$ws.Range("G12").Value = 19.63122766666666
$ws.Range("H12").Value = 58.893683
$ws.Range("I12").Value = 0.4223271270713038
$ws.Range("J12").Value = 0.4223271270713038
$ws.Range("M12").Value = 14.774121
$ws.Range("N12").Value = 44.322363
$ws.Range("O12").Value = 0.3559054518837227
$ws.Range("P12").Value = 0.3559054518837226
$ws.Range("Q12").Value = 290.034132925881
$ws.Range("R12").Value = 2610.307196332928
$ws.Range("S12").Value = 0.1503085270030667
$ws.Range("T12").Value = 0.1503085270030667

# Row 13
$ws.Range("G13").Value = 19.63122766666666
$ws.Range("H13").Value = 58.893683
$ws.Range("I13").Value = 0.4223271270713038
$ws.Range("J13").Value = 0.4223271270713038
$ws.Range("M13").Value = 16.53203833333334
$ws.Range("N13").Value = 49.596115
$ws.Range("O13").Value = 0.3982533088489005
$ws.Range("P13").Value = 0.3982533088489005
$ws.Range("Q13").Value = 324.5442083157272
$ws.Range("R13").Value = 2920.897874841545
$ws.Range("S13").Value = 0.1681931757727968
$ws.Range("T13").Value = 0.1681931757727968

# Row 14
$ws.Range("G14").Value = 6.251301333333333
$ws.Range("H14").Value = 18.753904
$ws.Range("I14").Value = 0.1344844131702722
$ws.Range("J14").Value = 0.1344844131702721
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.270625
$ws.Range("N14").Value = 0.811875
$ws.Range("O14").Value = 0.00651929904432436
$ws.Range("P14").Value = 0.006519299044324359
$ws.Range("Q14").Value = 1.691758423333333
$ws.Range("R14").Value = 15.22582581
$ws.Range("S14").Value = 0.0008767441062574777
$ws.Range("T14").Value = 0.0008767441062574774

# Row 15
$ws.Range("G15").Value = 6.251301333333333
$ws.Range("H15").Value = 18.753904
$ws.Range("I15").Value = 0.1344844131702722
$ws.Range("J15").Value = 0.1344844131702721
$ws.Range("O15").Value = 0.2393219402230525
$ws.Range("P15").Value = 0.2393219402230525
$ws.Range("Q15").Value = 62.10405528387378
$ws.Range("R15").Value = 558.936497554864
$ws.Range("S15").Value = 0.03218507068966817
$ws.Range("T15").Value = 0.03218507068966816

# Row 16
$ws.Range("G16").Value = 6.251301333333333
$ws.Range("H16").Value = 18.753904
$ws.Range("I16").Value = 0.1344844131702722
$ws.Range("J16").Value = 0.1344844131702721
$ws.Range("M16").Value = 14.774121
$ws.Range("N16").Value = 44.322363
$ws.Range("O16").Value = 0.3559054518837227
$ws.Range("P16").Value = 0.3559054518837226
$ws.Range("Q16").Value = 92.35748230612799
$ws.Range("R16").Value = 831.2173407551519
$ws.Range("S16").Value = 0.04786373584068298
$ws.Range("T16").Value = 0.04786373584068296

# Row 17
$ws.Range("G17").Value = 6.251301333333333
$ws.Range("H17").Value = 18.753904
$ws.Range("I17").Value = 0.1344844131702722
$ws.Range("J17").Value = 0.1344844131702721
$ws.Range("M17").Value = 16.53203833333334
$ws.Range("N17").Value = 49.596115
$ws.Range("O17").Value = 0.3982533088489005
$ws.Range("P17").Value = 0.3982533088489005
$ws.Range("Q17").Value = 103.3467532758845
$ws.Range("R17").Value = 930.12077948296
$ws.Range("S17").Value = 0.05355886253366354
$ws.Range("T17").Value = 0.05355886253366352
